$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "0.616") need to be
# forced to Text format first, otherwise Excel auto-converts them to a number
# (the source data is a text-formatted price column, e.g. "57.118.17").
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '57.375.19'
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").Value = '2.423.18'
$ws.Range("E3").Value = '  -2.51%  '

$ws.Range("E4").Value = '  -0.20%  '

Set-TextValue $ws.Range("D5") '488.97'
$ws.Range("E5").Value = '  -0.75%  '

Set-TextValue $ws.Range("D6") '153.88'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range("D7") '0.616'
$ws.Range("E7").Value = '  +20.01%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range("D8") '0.996'
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '2.443.73'
$ws.Range("E9").Value = '  -2.16%  '

Set-TextValue $ws.Range("D10") '6.18'
$ws.Range("E10").Value = '  +7.80%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("E13").Value = '  +1.44%  '

$ws.Range("D14").Value = '2.841.48'
$ws.Range("E14").Value = '  -3.04%  '

$ws.Range("D15").Value = '57.306.39'
$ws.Range("E15").Value = '  +0.17%  '

Set-TextValue $ws.Range("D16") '20.62'
$ws.Range("E16").Value = '  -2.47%  '

$ws.Range("E17").Value = '  -2.81%  '

$ws.Range("D18").Value = '2.435.57'
$ws.Range("E18").Value = '  -3.58%  '

Set-TextValue $ws.Range("D19") '4.66'
$ws.Range("E19").Value = '  +2.21%  '

Set-TextValue $ws.Range("D20") '324.82'
$ws.Range("E20").Value = '  +1.08%  '

Set-TextValue $ws.Range("D21") '10.03'
$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("E22").Value = '  -0.02%  '

Set-TextValue $ws.Range("D23") '5.96'
$ws.Range("E23").Value = '  +0.68%  '

Set-TextValue $ws.Range("D24") '57.95'
$ws.Range("E24").Value = '  -0.54%  '

$ws.Range("E25").Value = '  -1.30%  '

Set-TextValue $ws.Range("D26") '0.998'
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  -2.45%  '

$ws.Range("D28").Value = '2.529.32'
$ws.Range("E28").Value = '  -3.51%  '

$ws.Range("E29").Value = '  -3.81%  '

$ws.Range("E30").Value = '  -3.79%  '

Set-TextValue $ws.Range("D31") '0.999'
$ws.Range("E31").Value = '  -0.01%  '

Set-TextValue $ws.Range("D32") '151.47'
$ws.Range("E32").Value = '  +0.01%  '

Set-TextValue $ws.Range("D33") '18.65'
$ws.Range("E33").Value = '  +1.73%  '

$ws.Range("E34").Value = '  -0.49%  '

$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("E36").Value = '  -0.26%  '

$ws.Range("E37").Value = '  -0.47%  '

Set-TextValue $ws.Range("D38") '0.820'
$ws.Range("E38").Value = '  -9.11%  '

Set-TextValue $ws.Range("D39") '287.36'
$ws.Range("E39").Value = '  +9.25%  '

Set-TextValue $ws.Range("D40") '0.102'
$ws.Range("E40").Value = '  +7.82%  '

Set-TextValue $ws.Range("D41") '33.98'
$ws.Range("E41").Value = '  -0.83%  '

Set-TextValue $ws.Range("D42") '1.37'
$ws.Range("E42").Value = '  -1.82%  '

$ws.Range("E43").Value = '  +0.15%  '

Set-TextValue $ws.Range("D44") '0.996'
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("E45").Value = '  -2.26%  '

Set-TextValue $ws.Range("D46") '0.0531'
$ws.Range("E46").Value = '  -5.00%  '

Set-TextValue $ws.Range("D47") '10.21'
$ws.Range("E47").Value = '  -0.16%  '

Set-TextValue $ws.Range("D48") '0.0228'
$ws.Range("E48").Value = '  -0.44%  '

Set-TextValue $ws.Range("D49") '4.53'
$ws.Range("E49").Value = '  -6.96%  '

$ws.Range("D50").Value = '1.907.94'
$ws.Range("E50").Value = '  +1.10%  '

Set-TextValue $ws.Range("D51") '17.61'
$ws.Range("E51").Value = '  -1.74%  '

Write-Host "Applied crypto list update."